# Automatische sync: 2025-06-17 10:57:36
$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append new row (row 5) ----
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A5").Value = "Wat zijn jullie openingstijden?"
$wsLogs.Range("B5").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C5").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$wsLogs.Range("D5").Value = "Informatieaanvraag"
$wsLogs.Range("E5").Value = "Beste,`nBedankt voor je interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 uur tot 18:00 uur. Op zaterdag zijn wij geopend van 10:00 uur tot 16:00 uur. Voor eventuele feestdagen en afwijkende openingstijden, adviseer ik onze website te raadplegen.`nMet vriendelijke groet, [Jouw naam]"
$wsLogs.Range("F5").Value = "2025-06-17 10:28:27"
$wsLogs.Range("G5").Value = "Ja"

# Extend the conditional formatting ranges to include the new row
$wsLogs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D5"))
$wsLogs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G5"))

# ---- Dashboard sheet: append new row (row 5) ----
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A5").Value = "Informatieaanvraag"
$wsDash.Range("B5").Value = 1

# ---- Chart: extend category/value series ranges to row 5 ----
$chart = $wsDash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"
